$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 24V test data: the expected 24V PSU load values for the
# "addition and deletion of accessories" test row (row 8) were lowered.
$ws.Range("F8").Value = 0.329
$ws.Range("J8").Value = 0.381
$ws.Range("K8").Value = 0.329

# Reflect the author's saved selection/view state on the sheet.
$ws.Activate()
$ws.Range("J19").Select()
